$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.275.87"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "3.131.89"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'603.96"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").Value = "'143.11"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.129.17"
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("D9").Value = "'0.523"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("E11").Value = "  +2.95%  "
$ws.Range("D12").Value = "'0.469"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("E13").Value = "  +3.23%  "
$ws.Range("D14").Value = "'35.19"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "3.641.53"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").Value = "64.184.60"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "3.126.94"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "'6.90"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").Value = "'480.20"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").Value = "'14.58"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").Value = "'7.70"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "'85.55"
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("D25").Value = "'13.43"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "'8.38"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +7.46%  "
$ws.Range("D30").Value = "'2.05"
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "'26.90"
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("D34").Value = "'2.65"
$ws.Range("E34").Value = "  -2.25%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0772"
$ws.Range("E36").Value = "  +5.15%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "'5.97"
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("D38").Value = "'52.39"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").Value = "'3.02"
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("D40").Value = "'445.92"
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("D41").Value = "'0.0393"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").Value = "'8.21"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "2.852.88"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").Value = "'0.261"
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").Value = "'26.08"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").Value = "'120.45"
$ws.Range("E51").Value = "  +1.87%  "
